# Updates the cryptos price/volume snapshot in the active worksheet.
# Each entry: cell reference, new text value, and whether the value must be
# forced into Excel as literal text (leading apostrophe) because it would
# otherwise be auto-parsed as a number (e.g. "570.12" -> 570.12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "64.259.64"; ForceText = $false }
    @{ Cell = "E2"; Value = "  -1.30%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "3.398.84"; ForceText = $false }
    @{ Cell = "E3"; Value = "  -1.38%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  +0.02%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "570.12"; ForceText = $true }
    @{ Cell = "E5"; Value = "  -1.05%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "156.01"; ForceText = $true }
    @{ Cell = "E6"; Value = "  -3.26%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "0.630"; ForceText = $true }
    @{ Cell = "E7"; Value = "  +8.74%  "; ForceText = $false }
    @{ Cell = "E8"; Value = "  +0.00%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "3.402.02"; ForceText = $false }
    @{ Cell = "E9"; Value = "  -1.42%  "; ForceText = $false }
    @{ Cell = "E10"; Value = "  -1.97%  "; ForceText = $false }
    @{ Cell = "E11"; Value = "  -2.14%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "0.438"; ForceText = $true }
    @{ Cell = "E12"; Value = "  +0.15%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "3.987.33"; ForceText = $false }
    @{ Cell = "E13"; Value = "  -1.33%  "; ForceText = $false }
    @{ Cell = "E14"; Value = "  -0.25%  "; ForceText = $false }
    @{ Cell = "E15"; Value = "  -3.34%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "27.36"; ForceText = $true }
    @{ Cell = "E16"; Value = "  -2.60%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "64.283.51"; ForceText = $false }
    @{ Cell = "E17"; Value = "  -1.15%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "3.372.45"; ForceText = $false }
    @{ Cell = "E18"; Value = "  -3.83%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "6.28"; ForceText = $true }
    @{ Cell = "E19"; Value = "  -1.25%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "13.79"; ForceText = $true }
    @{ Cell = "E20"; Value = "  -3.25%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "377.39"; ForceText = $true }
    @{ Cell = "E21"; Value = "  -2.62%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "7.97"; ForceText = $true }
    @{ Cell = "E22"; Value = "  -2.70%  "; ForceText = $false }
    @{ Cell = "B23"; Value = "Dai"; ForceText = $false }
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; ForceText = $false }
    @{ Cell = "D23"; Value = "1.00"; ForceText = $true }
    @{ Cell = "E23"; Value = "  +0.15%  "; ForceText = $false }
    @{ Cell = "B24"; Value = "Polygon"; ForceText = $false }
    @{ Cell = "C24"; Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; ForceText = $false }
    @{ Cell = "D24"; Value = "0.541"; ForceText = $true }
    @{ Cell = "E24"; Value = "  -0.38%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "71.66"; ForceText = $true }
    @{ Cell = "E25"; Value = "  -1.92%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "0.0000118"; ForceText = $true }
    @{ Cell = "E26"; Value = "  -5.02%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "10.35"; ForceText = $true }
    @{ Cell = "E27"; Value = "  +7.31%  "; ForceText = $false }
    @{ Cell = "E28"; Value = "  -2.01%  "; ForceText = $false }
    @{ Cell = "E29"; Value = "  -0.03%  "; ForceText = $false }
    @{ Cell = "E30"; Value = "  +2.03%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "6.09"; ForceText = $true }
    @{ Cell = "E31"; Value = "  -2.53%  "; ForceText = $false }
    @{ Cell = "E32"; Value = "  -1.77%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "23.01"; ForceText = $true }
    @{ Cell = "E33"; Value = "  -2.70%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "7.08"; ForceText = $true }
    @{ Cell = "E34"; Value = "  +0.33%  "; ForceText = $false }
    @{ Cell = "E35"; Value = "  +6.08%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "160.62"; ForceText = $true }
    @{ Cell = "E36"; Value = "  -0.63%  "; ForceText = $false }
    @{ Cell = "E37"; Value = "  -1.62%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "0.0757"; ForceText = $true }
    @{ Cell = "E38"; Value = "  -1.26%  "; ForceText = $false }
    @{ Cell = "B39"; Value = "Maker"; ForceText = $false }
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; ForceText = $false }
    @{ Cell = "D39"; Value = "2.867.02"; ForceText = $false }
    @{ Cell = "E39"; Value = "  -5.61%  "; ForceText = $false }
    @{ Cell = "B40"; Value = "RenderToken"; ForceText = $false }
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; ForceText = $false }
    @{ Cell = "D40"; Value = "6.73"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +1.78%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "26.31"; ForceText = $true }
    @{ Cell = "E41"; Value = "  -3.47%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "4.57"; ForceText = $true }
    @{ Cell = "E42"; Value = "  +0.73%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "42.63"; ForceText = $true }
    @{ Cell = "E43"; Value = "  -0.25%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "0.0313"; ForceText = $true }
    @{ Cell = "E44"; Value = "  -0.91%  "; ForceText = $false }
    @{ Cell = "B45"; Value = "Mantle"; ForceText = $false }
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; ForceText = $false }
    @{ Cell = "D45"; Value = "0.767"; ForceText = $true }
    @{ Cell = "E45"; Value = "  -0.51%  "; ForceText = $false }
    @{ Cell = "B46"; Value = "InjectiveProtocol"; ForceText = $false }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; ForceText = $false }
    @{ Cell = "D46"; Value = "25.73"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +3.46%  "; ForceText = $false }
    @{ Cell = "D47"; Value = "321.10"; ForceText = $true }
    @{ Cell = "E47"; Value = "  +5.23%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "0.110"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +3.48%  "; ForceText = $false }
    @{ Cell = "E49"; Value = "  -2.11%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "2.18"; ForceText = $true }
    @{ Cell = "E50"; Value = "  -0.95%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "6.51"; ForceText = $true }
    @{ Cell = "E51"; Value = "  -1.48%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Leading apostrophe tells Excel to store the value as plain text
        # instead of silently converting it to a number.
        $rng.Value = "'" + $u.Value
        # Reset style so the quote-prefix formatting Excel applies does not
        # leave a stray style/number-format difference on the cell.
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
